# Update pl_mw.xlsx results for the 380 kV case
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.214278702769434
$ws.Range("C2").Value = 0.365622547721614
$ws.Range("D2").Value = 0.07941296282635335
$ws.Range("E2").Value = 0.4223042686625433
$ws.Range("G2").Value = 0.002386932758169901
$ws.Range("I2").Value = 0.3355830042883845
$ws.Range("O2").Value = 1.839492726063412
$ws.Range("B3").Value = 1.066785532073254
$ws.Range("C3").Value = 0.3189329536006085
$ws.Range("D3").Value = 0.07190151070420825
$ws.Range("E3").Value = 0.3682411441012619
$ws.Range("G3").Value = 0.0023906544224222
$ws.Range("I3").Value = 0.342846019663682
$ws.Range("O3").Value = 1.829181641849857
$ws.Range("B4").Value = 0.9760460895228107
$ws.Range("C4").Value = 0.2901710387618266
$ws.Range("D4").Value = 0.06732703250825978
$ws.Range("E4").Value = 0.3351553121666768
$ws.Range("G4").Value = 0.0023930571205127
$ws.Range("I4").Value = 0.3477362606604721
$ws.Range("O4").Value = 1.825035712061833
$ws.Range("B5").Value = 0.939025595801354
$ws.Range("C5").Value = 0.2784268330206032
$ws.Range("D5").Value = 0.06547225941756096
$ws.Range("E5").Value = 0.3216977307559574
$ws.Range("G5").Value = 0.002394065907033896
$ws.Range("I5").Value = 0.3498367422637827
$ws.Range("O5").Value = 1.823891709837852
$ws.Range("B6").Value = 0.932875787872149
$ws.Range("C6").Value = 0.2764753046428723
$ws.Range("D6").Value = 0.06516483880686508
$ws.Range("E6").Value = 0.3194645631544688
$ws.Range("G6").Value = 0.002394235210051644
$ws.Range("I6").Value = 0.3501920098564959
$ws.Range("O6").Value = 1.823734579946034
$ws.Range("B7").Value = 0.9755469922606608
$ws.Range("C7").Value = 0.2900127471027929
$ws.Range("D7").Value = 0.0673019805883257
$ws.Range("E7").Value = 0.3349737196120941
$ws.Range("G7").Value = 0.002393070605070315
$ws.Range("I7").Value = 0.3477641533952962
$ws.Range("O7").Value = 1.825018080220559
$ws.Range("B8").Value = 1.163460754320909
$ws.Range("C8").Value = 0.3495436679585282
$ws.Range("D8").Value = 0.07681515872241107
$ws.Range("E8").Value = 0.4036392564278231
$ws.Range("G8").Value = 0.002388191644845439
$ws.Range("I8").Value = 0.3379975116692222
$ws.Range("O8").Value = 1.835481543054357
$ws.Range("B9").Value = 1.530507497587621
$ws.Range("C9").Value = 0.4655316202233735
$ws.Range("D9").Value = 0.09577362304975168
$ws.Range("E9").Value = 0.5392677428273771
$ws.Range("G9").Value = 0.002379552355474105
$ws.Range("I9").Value = 0.3222913369731444
$ws.Range("O9").Value = 1.873521689936382
$ws.Range("B10").Value = 1.799268722748423
$ws.Range("C10").Value = 0.5502938090637599
$ws.Range("D10").Value = 0.1098955321608202
$ws.Range("E10").Value = 0.6396700428514492
$ws.Range("G10").Value = 0.002373764517215854
$ws.Range("I10").Value = 0.3128921671343576
$ws.Range("O10").Value = 1.912411571984819
$ws.Range("B11").Value = 1.921335567927997
$ws.Range("C11").Value = 0.5887573983310972
$ws.Range("D11").Value = 0.1163637848113126
$ws.Range("E11").Value = 0.6855455647639701
$ws.Range("G11").Value = 0.002371251569842677
$ws.Range("I11").Value = 0.3090894342826473
$ws.Range("O11").Value = 1.932536207612713
$ws.Range("B12").Value = 1.967530449814149
$ws.Range("C12").Value = 0.6033088396283119
$ws.Range("D12").Value = 0.1188196122638772
$ws.Range("E12").Value = 0.7029492356389682
$ws.Range("G12").Value = 0.002370317127808108
$ws.Range("I12").Value = 0.3077181603744634
$ws.Range("O12").Value = 1.9405113263484
$ws.Range("B13").Value = 1.957582873452964
$ws.Range("C13").Value = 0.6001755470046533
$ws.Range("D13").Value = 0.1182904179037649
$ws.Range("E13").Value = 0.6991995969189873
$ws.Range("G13").Value = 0.002370517614953276
$ws.Range("I13").Value = 0.3080104200446456
$ws.Range("O13").Value = 1.93877791091316
$ws.Range("B14").Value = 1.925136640763355
$ws.Range("C14").Value = 0.589954833552099
$ws.Range("D14").Value = 0.1165656975841785
$ws.Range("E14").Value = 0.6869767270220137
$ws.Range("G14").Value = 0.002371174349446189
$ws.Range("I14").Value = 0.3089752369824446
$ws.Range("O14").Value = 1.933185198136187
$ws.Range("B15").Value = 1.905258548358972
$ws.Range("C15").Value = 0.5836925374246675
$ws.Range("D15").Value = 0.1155100981930417
$ws.Range("E15").Value = 0.6794940588900573
$ws.Range("G15").Value = 0.00237157884972322
$ws.Range("I15").Value = 0.309575188487976
$ws.Range("O15").Value = 1.929805779635217
$ws.Range("B16").Value = 1.791287352754466
$ws.Range("C16").Value = 0.5477781907008534
$ws.Range("D16").Value = 0.1094737126759497
$ws.Range("E16").Value = 0.6366762523799139
$ws.Range("G16").Value = 0.002373931150105051
$ws.Range("I16").Value = 0.3131502642231609
$ws.Range("O16").Value = 1.911145724827605
$ws.Range("B17").Value = 1.721319101335837
$ws.Range("C17").Value = 0.5257213826342877
$ws.Range("D17").Value = 0.1057819546881404
$ws.Range("E17").Value = 0.6104624632130822
$ws.Range("G17").Value = 0.002375404869659113
$ws.Range("I17").Value = 0.3154651158991157
$ws.Range("O17").Value = 1.900324791078134
$ws.Range("B18").Value = 1.681057029281135
$ws.Range("C18").Value = 0.5130259447262233
$ws.Range("D18").Value = 0.1036627006402142
$ws.Range("E18").Value = 0.5954038442336156
$ws.Range("G18").Value = 0.002376263811799212
$ws.Range("I18").Value = 0.3168410346739954
$ws.Range("O18").Value = 1.894329656059483
$ws.Range("B19").Value = 1.667421900263832
$ws.Range("C19").Value = 0.5087259576068845
$ws.Range("D19").Value = 0.102945867087314
$ws.Range("E19").Value = 0.5903084163634844
$ws.Range("G19").Value = 0.002376556578046494
$ws.Range("I19").Value = 0.3173145141980172
$ws.Range("O19").Value = 1.892338964883692
$ws.Range("B20").Value = 1.72876923261822
$ws.Range("C20").Value = 0.5280702948022622
$ws.Range("D20").Value = 0.1061745186881922
$ws.Range("E20").Value = 0.6132509955064762
$ws.Range("G20").Value = 0.002375246821049981
$ws.Range("I20").Value = 0.3152140875479219
$ws.Range("O20").Value = 1.901452988148264
$ws.Range("B21").Value = 1.934667687890965
$ws.Range("C21").Value = 0.5929572832286567
$ws.Range("D21").Value = 0.1170721141883462
$ws.Range("E21").Value = 0.6905660018494189
$ws.Range("G21").Value = 0.002370980985845588
$ws.Range("I21").Value = 0.3086899755443433
$ws.Range("O21").Value = 1.934818260883503
$ws.Range("B22").Value = 2.069063561715268
$ws.Range("C22").Value = 0.6352836636619372
$ws.Range("D22").Value = 0.124231908524834
$ws.Range("E22").Value = 0.7412812670968663
$ws.Range("G22").Value = 0.002368292972933228
$ws.Range("I22").Value = 0.3048270894215861
$ws.Range("O22").Value = 1.958691928414339
$ws.Range("B23").Value = 1.997349985332391
$ws.Range("C23").Value = 0.6127007541232956
$ws.Range("D23").Value = 0.1204071205783634
$ws.Range("E23").Value = 0.7141957533454217
$ws.Range("G23").Value = 0.002369718501079297
$ws.Range("I23").Value = 0.3068518535424296
$ws.Range("O23").Value = 1.945759452377757
$ws.Range("B24").Value = 1.725401141527641
$ws.Range("C24").Value = 0.5270083973428541
$ws.Range("D24").Value = 0.1059970305859821
$ws.Range("E24").Value = 0.6119902627804805
$ws.Range("G24").Value = 0.002375318238429122
$ws.Range("I24").Value = 0.3153274370996826
$ws.Range("O24").Value = 1.900942226728915
$ws.Range("B25").Value = 1.431369578700071
$ws.Range("C25").Value = 0.4342339205458075
$ws.Range("D25").Value = 0.09061154676106753
$ws.Range("E25").Value = 0.5024556456920095
$ws.Range("G25").Value = 0.002381790799879229
$ws.Range("I25").Value = 0.326167287838647
$ws.Range("O25").Value = 1.861328217133774
